$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Add formulas to the numbering column (A21:A29), building a running count
$ws.Range("A21").Formula = "=A20+1"
$ws.Range("A22:A29").Formula = "=A21+1"

# Rows 28/29 previously had the "last row" border style (no bottom hairline);
# now that row 29 gets data too, copy the interior-row format from A20.
$ws.Range("A20").Copy()
$ws.Range("A28:A29").PasteSpecial(-4122)

# New row 29 data: QUERY_TIMEOUT_DEFAULT constant
$ws.Range("B29").Value = "QUERY_TIMEOUT_DEFAULT"
$ws.Range("C29").Value = "java.lang.Long"
$ws.Range("D29").Value = "60000L"
$ws.Range("E29").Value = "クエリタイムアウトを強制指定する場合のデフォルト値（ミリ秒）です。"

# Update selection to match the final state
[void]$ws.Range("E33").Select()
